$wb = $excel.ActiveWorkbook

# Rename the sole worksheet from "AppUserStore" to the Vietnamese title
# "Phạm vi đi tuyến" (DMS: Translate AppUserStoreMapping Export and ExportTemplate)
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Phạm vi đi tuyến"
